$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 2.338576666666667
$ws.Cells.Item(2, 8).Value = 7.01573
$ws.Cells.Item(2, 9).Value = 0.06036918779536898
$ws.Cells.Item(2, 10).Value = 0.06036918779536898
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 15.85269333333333
$ws.Cells.Item(2, 14).Value = 47.55808
$ws.Cells.Item(2, 15).Value = 0.363669005908797
$ws.Cells.Item(2, 16).Value = 0.363669005908797
$ws.Cells.Item(2, 17).Value = 37.07273873315556
$ws.Cells.Item(2, 18).Value = 333.6546485984
$ws.Cells.Item(2, 19).Value = 0.02195440251306332
$ws.Cells.Item(2, 20).Value = 0.02195440251306332

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 2.338576666666667
$ws.Cells.Item(3, 8).Value = 7.01573
$ws.Cells.Item(3, 9).Value = 0.06036918779536898
$ws.Cells.Item(3, 10).Value = 0.06036918779536898
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 9.934580333333335
$ws.Cells.Item(3, 14).Value = 29.803741
$ws.Cells.Item(3, 15).Value = 0.2279044246915194
$ws.Cells.Item(3, 16).Value = 0.2279044246915194
$ws.Cells.Item(3, 17).Value = 23.23277776065889
$ws.Cells.Item(3, 18).Value = 209.09499984593
$ws.Cells.Item(3, 19).Value = 0.01375840501359787
$ws.Cells.Item(3, 20).Value = 0.01375840501359787

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 2.338576666666667
$ws.Cells.Item(4, 8).Value = 7.01573
$ws.Cells.Item(4, 9).Value = 0.06036918779536898
$ws.Cells.Item(4, 10).Value = 0.06036918779536898
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 2.849702333333333
$ws.Cells.Item(4, 14).Value = 8.549106999999999
$ws.Cells.Item(4, 15).Value = 0.06537364931674991
$ws.Cells.Item(4, 16).Value = 0.06537364931674992
$ws.Cells.Item(4, 17).Value = 6.664247383678889
$ws.Cells.Item(4, 18).Value = 59.97822645311
$ws.Cells.Item(4, 19).Value = 0.00394655411247147
$ws.Cells.Item(4, 20).Value = 0.003946554112471471

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 2.338576666666667
$ws.Cells.Item(5, 8).Value = 7.01573
$ws.Cells.Item(5, 9).Value = 0.06036918779536898
$ws.Cells.Item(5, 10).Value = 0.06036918779536898
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 14.95401766666667
$ws.Cells.Item(5, 14).Value = 44.862053
$ws.Cells.Item(5, 15).Value = 0.3430529200829336
$ws.Cells.Item(5, 16).Value = 0.3430529200829337
$ws.Cells.Item(5, 17).Value = 34.97111678818779
$ws.Cells.Item(5, 18).Value = 314.7400510936901
$ws.Cells.Item(5, 19).Value = 0.02070982615623633
$ws.Cells.Item(5, 20).Value = 0.02070982615623633

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 15.98183966666667
$ws.Cells.Item(6, 8).Value = 47.945519
$ws.Cells.Item(6, 9).Value = 0.4125632030390895
$ws.Cells.Item(6, 10).Value = 0.4125632030390896
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 15.85269333333333
$ws.Cells.Item(6, 14).Value = 47.55808
$ws.Cells.Item(6, 15).Value = 0.363669005908797
$ws.Cells.Item(6, 16).Value = 0.363669005908797
$ws.Cells.Item(6, 17).Value = 253.3552031381689
$ws.Cells.Item(6, 18).Value = 2280.19682824352
$ws.Cells.Item(6, 19).Value = 0.1500364499237749
$ws.Cells.Item(6, 20).Value = 0.1500364499237749

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 15.98183966666667
$ws.Cells.Item(7, 8).Value = 47.945519
$ws.Cells.Item(7, 9).Value = 0.4125632030390895
$ws.Cells.Item(7, 10).Value = 0.4125632030390896
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 9.934580333333335
$ws.Cells.Item(7, 14).Value = 29.803741
$ws.Cells.Item(7, 15).Value = 0.2279044246915194
$ws.Cells.Item(7, 16).Value = 0.2279044246915194
$ws.Cells.Item(7, 17).Value = 158.7728700429533
$ws.Cells.Item(7, 18).Value = 1428.955830386579
$ws.Cells.Item(7, 19).Value = 0.09402497943751421
$ws.Cells.Item(7, 20).Value = 0.09402497943751423

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 15.98183966666667
$ws.Cells.Item(8, 8).Value = 47.945519
$ws.Cells.Item(8, 9).Value = 0.4125632030390895
$ws.Cells.Item(8, 10).Value = 0.4125632030390896
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 2.849702333333333
$ws.Cells.Item(8, 14).Value = 8.549106999999999
$ws.Cells.Item(8, 15).Value = 0.06537364931674991
$ws.Cells.Item(8, 16).Value = 0.06537364931674992
$ws.Cells.Item(8, 17).Value = 45.54348578905922
$ws.Cells.Item(8, 18).Value = 409.891372101533
$ws.Cells.Item(8, 19).Value = 0.02697076215647253
$ws.Cells.Item(8, 20).Value = 0.02697076215647254

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 15.98183966666667
$ws.Cells.Item(9, 8).Value = 47.945519
$ws.Cells.Item(9, 9).Value = 0.4125632030390895
$ws.Cells.Item(9, 10).Value = 0.4125632030390896
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 14.95401766666667
$ws.Cells.Item(9, 14).Value = 44.862053
$ws.Cells.Item(9, 15).Value = 0.3430529200829336
$ws.Cells.Item(9, 16).Value = 0.3430529200829337
$ws.Cells.Item(9, 17).Value = 238.9927127211675
$ws.Cells.Item(9, 18).Value = 2150.934414490507
$ws.Cells.Item(9, 19).Value = 0.1415310115213279
$ws.Cells.Item(9, 20).Value = 0.1415310115213279

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 4.530968000000001
$ws.Cells.Item(10, 8).Value = 13.592904
$ws.Cells.Item(10, 9).Value = 0.1169646742762937
$ws.Cells.Item(10, 10).Value = 0.1169646742762937
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 15.85269333333333
$ws.Cells.Item(10, 14).Value = 47.55808
$ws.Cells.Item(10, 15).Value = 0.363669005908797
$ws.Cells.Item(10, 16).Value = 0.363669005908797
$ws.Cells.Item(10, 17).Value = 71.82804620714667
$ws.Cells.Item(10, 18).Value = 646.45241586432
$ws.Cells.Item(10, 19).Value = 0.04253642682050598
$ws.Cells.Item(10, 20).Value = 0.04253642682050598

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 4.530968000000001
$ws.Cells.Item(11, 8).Value = 13.592904
$ws.Cells.Item(11, 9).Value = 0.1169646742762937
$ws.Cells.Item(11, 10).Value = 0.1169646742762937
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 9.934580333333335
$ws.Cells.Item(11, 14).Value = 29.803741
$ws.Cells.Item(11, 15).Value = 0.2279044246915194
$ws.Cells.Item(11, 16).Value = 0.2279044246915194
$ws.Cells.Item(11, 17).Value = 45.01326558376268
$ws.Cells.Item(11, 18).Value = 405.119390253864
$ws.Cells.Item(11, 19).Value = 0.02665676680016969
$ws.Cells.Item(11, 20).Value = 0.02665676680016969

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 4.530968000000001
$ws.Cells.Item(12, 8).Value = 13.592904
$ws.Cells.Item(12, 9).Value = 0.1169646742762937
$ws.Cells.Item(12, 10).Value = 0.1169646742762937
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 2.849702333333333
$ws.Cells.Item(12, 14).Value = 8.549106999999999
$ws.Cells.Item(12, 15).Value = 0.06537364931674991
$ws.Cells.Item(12, 16).Value = 0.06537364931674992
$ws.Cells.Item(12, 17).Value = 12.91191008185867
$ws.Cells.Item(12, 18).Value = 116.207190736728
$ws.Cells.Item(12, 19).Value = 0.007646407598586306
$ws.Cells.Item(12, 20).Value = 0.007646407598586307

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 4.530968000000001
$ws.Cells.Item(13, 8).Value = 13.592904
$ws.Cells.Item(13, 9).Value = 0.1169646742762937
$ws.Cells.Item(13, 10).Value = 0.1169646742762937
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 14.95401766666667
$ws.Cells.Item(13, 14).Value = 44.862053
$ws.Cells.Item(13, 15).Value = 0.3430529200829336
$ws.Cells.Item(13, 16).Value = 0.3430529200829337
$ws.Cells.Item(13, 17).Value = 67.75617551910135
$ws.Cells.Item(13, 18).Value = 609.805579671912
$ws.Cells.Item(13, 19).Value = 0.04012507305703176
$ws.Cells.Item(13, 20).Value = 0.04012507305703177

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 15.886534
$ws.Cells.Item(14, 8).Value = 47.659602
$ws.Cells.Item(14, 9).Value = 0.4101029348892478
$ws.Cells.Item(14, 10).Value = 0.4101029348892479
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 15.85269333333333
$ws.Cells.Item(14, 14).Value = 47.55808
$ws.Cells.Item(14, 15).Value = 0.363669005908797
$ws.Cells.Item(14, 16).Value = 0.363669005908797
$ws.Cells.Item(14, 17).Value = 251.8443516315733
$ws.Cells.Item(14, 18).Value = 2266.59916468416
$ws.Cells.Item(14, 19).Value = 0.1491417266514528
$ws.Cells.Item(14, 20).Value = 0.1491417266514529

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 15.886534
$ws.Cells.Item(15, 8).Value = 47.659602
$ws.Cells.Item(15, 9).Value = 0.4101029348892478
$ws.Cells.Item(15, 10).Value = 0.4101029348892479
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 9.934580333333335
$ws.Cells.Item(15, 14).Value = 29.803741
$ws.Cells.Item(15, 15).Value = 0.2279044246915194
$ws.Cells.Item(15, 16).Value = 0.2279044246915194
$ws.Cells.Item(15, 17).Value = 157.8260482412313
$ws.Cells.Item(15, 18).Value = 1420.434434171082
$ws.Cells.Item(15, 19).Value = 0.09346427344023768
$ws.Cells.Item(15, 20).Value = 0.09346427344023769

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 15.886534
$ws.Cells.Item(16, 8).Value = 47.659602
$ws.Cells.Item(16, 9).Value = 0.4101029348892478
$ws.Cells.Item(16, 10).Value = 0.4101029348892479
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 2.849702333333333
$ws.Cells.Item(16, 14).Value = 8.549106999999999
$ws.Cells.Item(16, 15).Value = 0.06537364931674991
$ws.Cells.Item(16, 16).Value = 0.06537364931674992
$ws.Cells.Item(16, 17).Value = 45.27189300837933
$ws.Cells.Item(16, 18).Value = 407.447037075414
$ws.Cells.Item(16, 19).Value = 0.02680992544921961
$ws.Cells.Item(16, 20).Value = 0.02680992544921962

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 15.886534
$ws.Cells.Item(17, 8).Value = 47.659602
$ws.Cells.Item(17, 9).Value = 0.4101029348892478
$ws.Cells.Item(17, 10).Value = 0.4101029348892479
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 14.95401766666667
$ws.Cells.Item(17, 14).Value = 44.862053
$ws.Cells.Item(17, 15).Value = 0.3430529200829336
$ws.Cells.Item(17, 16).Value = 0.3430529200829337
$ws.Cells.Item(17, 17).Value = 237.5675100981007
$ws.Cells.Item(17, 18).Value = 2138.107590882906
$ws.Cells.Item(17, 19).Value = 0.1406870093483377
$ws.Cells.Item(17, 20).Value = 0.1406870093483377
